$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2491.5
$ws.Range("J17").Value = 2491.5
$ws.Range("L17").Value = 7474.5
$ws.Range("N17").Value = -7810.5
$ws.Range("H86").Value = 1727.4375
$ws.Range("I86").Value = 1661.4615
$ws.Range("J86").Value = 2013.3334
$ws.Range("K86").Value = 1661.4615
$ws.Range("L86").Value = 2013.3334
$ws.Range("M86").Value = -538.4614999999999
$ws.Range("N86").Value = -4259.3334
$ws.Range("H89").Value = 1727.4375
$ws.Range("I89").Value = 1661.4615
$ws.Range("J89").Value = 2013.3334
$ws.Range("K89").Value = 8307.307499999999
$ws.Range("L89").Value = 10066.667
$ws.Range("M89").Value = -2691.307499999999
$ws.Range("N89").Value = -21298.667
$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 5000
$ws.Range("N105").Value = -11988
$ws.Range("H137").Value = 28217.82
$ws.Range("I137").Value = 34519.902
$ws.Range("K137").Value = 103559.706
$ws.Range("M137").Value = -101009.706
$ws.Range("H138").Value = 2459.9375
$ws.Range("I138").Value = 2110.2334
$ws.Range("J138").Value = 3042.7778
$ws.Range("K138").Value = 6330.7002
$ws.Range("L138").Value = 9128.3334
$ws.Range("M138").Value = -1190.7002
$ws.Range("N138").Value = -19408.3334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2524.2964
$ws.Range("I61").Value = 2429.077
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2429.077
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2217.077
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 1971.75
$ws.Range("I74").Value = 1966.1
$ws.Range("K74").Value = 1966.1
$ws.Range("M74").Value = -1092.1
$ws.Range("H77").Value = 1971.75
$ws.Range("I77").Value = 1966.1
$ws.Range("K77").Value = 9830.5
$ws.Range("M77").Value = -5462.5
$ws.Range("H122").Value = 1996.0416
$ws.Range("I122").Value = 1924.1904
$ws.Range("J122").Value = 2499
$ws.Range("K122").Value = 5772.5712
$ws.Range("L122").Value = 7497
$ws.Range("M122").Value = -3322.5712
$ws.Range("N122").Value = -12397
$ws.Range("H132").Value = 29723.158
$ws.Range("I132").Value = 39314.035
$ws.Range("J132").Value = 2868.7
$ws.Range("K132").Value = 117942.105
$ws.Range("L132").Value = 8606.099999999999
$ws.Range("M132").Value = -115412.105
$ws.Range("N132").Value = -13666.1
$ws.Range("H136").Value = 2524.2964
$ws.Range("I136").Value = 2429.077
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7287.231000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4737.231000000001
$ws.Range("N136").Value = -20100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 28775.6
$ws.Range("J76").Value = 28775.6
$ws.Range("L76").Value = 28775.6
$ws.Range("N76").Value = -29405.6
$ws.Range("H79").Value = 28775.6
$ws.Range("J79").Value = 28775.6
$ws.Range("L79").Value = 28775.6
$ws.Range("N79").Value = -30959.6
$ws.Range("H86").Value = 3591.2144
$ws.Range("I86").Value = 3819.125
$ws.Range("J86").Value = 3287.3333
$ws.Range("K86").Value = 3819.125
$ws.Range("L86").Value = 3287.3333
$ws.Range("M86").Value = -2696.125
$ws.Range("N86").Value = -5533.3333
$ws.Range("H88").Value = 21541.857
$ws.Range("J88").Value = 21541.857
$ws.Range("L88").Value = 21541.857
$ws.Range("N88").Value = -22353.857
$ws.Range("H89").Value = 3591.2144
$ws.Range("I89").Value = 3819.125
$ws.Range("J89").Value = 3287.3333
$ws.Range("K89").Value = 19095.625
$ws.Range("L89").Value = 16436.6665
$ws.Range("M89").Value = -13479.625
$ws.Range("N89").Value = -27668.6665
$ws.Range("H91").Value = 21541.857
$ws.Range("J91").Value = 21541.857
$ws.Range("L91").Value = 21541.857
$ws.Range("N91").Value = -24349.857

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 165.94737
$ws.Range("I7").Value = 131.09091
$ws.Range("K7").Value = 131.09091
$ws.Range("M7").Value = -18.09091000000001
$ws.Range("H31").Value = 1863.6666
$ws.Range("I31").Value = 1143.6
$ws.Range("J31").Value = 3663.8333
$ws.Range("K31").Value = 1143.6
$ws.Range("L31").Value = 3663.8333
$ws.Range("M31").Value = -848.5999999999999
$ws.Range("N31").Value = -4253.8333
$ws.Range("H34").Value = 1863.6666
$ws.Range("I34").Value = 1143.6
$ws.Range("J34").Value = 3663.8333
$ws.Range("K34").Value = 1143.6
$ws.Range("L34").Value = 3663.8333
$ws.Range("M34").Value = -941.5999999999999
$ws.Range("N34").Value = -4067.8333
$ws.Range("H43").Value = 22029.416
$ws.Range("J43").Value = 22029.416
$ws.Range("L43").Value = 22029.416
$ws.Range("N43").Value = -22397.416
$ws.Range("H58").Value = 22326.334
$ws.Range("I58").Value = 24480.627
$ws.Range("J58").Value = 3799.4
$ws.Range("K58").Value = 24480.627
$ws.Range("L58").Value = 3799.4
$ws.Range("M58").Value = -24277.627
$ws.Range("N58").Value = -4205.4
$ws.Range("H88").Value = 19320.5
$ws.Range("J88").Value = 18651.857
$ws.Range("L88").Value = 18651.857
$ws.Range("N88").Value = -19463.857
$ws.Range("H91").Value = 19320.5
$ws.Range("J91").Value = 18651.857
$ws.Range("L91").Value = 18651.857
$ws.Range("N91").Value = -21459.857
$ws.Range("H99").Value = 7499.5
$ws.Range("J99").Value = 5999.3335
$ws.Range("L99").Value = 5999.3335
$ws.Range("N99").Value = -8995.3335
$ws.Range("H101").Value = 22029.416
$ws.Range("J101").Value = 22029.416
$ws.Range("L101").Value = 22029.416
$ws.Range("N101").Value = -28519.416
$ws.Range("H126").Value = 7499.5
$ws.Range("J126").Value = 5999.3335
$ws.Range("L126").Value = 17998.0005
$ws.Range("N126").Value = -22938.0005
$ws.Range("H132").Value = 2995.3635
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2995.3635
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 8986.0905
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14046.0905
$ws.Range("H136").Value = 22326.334
$ws.Range("I136").Value = 24480.627
$ws.Range("J136").Value = 3799.4
$ws.Range("K136").Value = 73441.881
$ws.Range("L136").Value = 11398.2
$ws.Range("M136").Value = -70891.881
$ws.Range("N136").Value = -16498.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2091.3333
$ws.Range("I136").Value = 1894.2106
$ws.Range("J136").Value = 2840.4
$ws.Range("K136").Value = 5682.6318
$ws.Range("L136").Value = 8521.2
$ws.Range("M136").Value = -3132.6318
$ws.Range("N136").Value = -13621.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 12022
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H96").Value = 2191.875
$ws.Range("I96").Value = 2112
$ws.Range("J96").Value = 2271.75
$ws.Range("K96").Value = 2112
$ws.Range("L96").Value = 2271.75
$ws.Range("M96").Value = -739
$ws.Range("N96").Value = -5017.75
$ws.Range("H100").Value = 492.7143
$ws.Range("I100").Value = 411.14285
$ws.Range("J100").Value = 655.8571
$ws.Range("K100").Value = 822.2857
$ws.Range("L100").Value = 1311.7142
$ws.Range("M100").Value = -281.2857
$ws.Range("N100").Value = -2393.7142
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
